$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

$ws.Range("J2").Value = 0.0104
$ws.Range("K2").Value = 0.2732
